$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 490.63635
$ws.Range("I15").Value = 490.63635
$ws.Range("K15").Value = 1471.90905
$ws.Range("M15").Value = -1302.90905
$ws.Range("H40").Value = 1900
$ws.Range("I40").Value = 1490
$ws.Range("J40").Value = 6000
$ws.Range("K40").Value = 1490
$ws.Range("L40").Value = 6000
$ws.Range("M40").Value = -1315
$ws.Range("N40").Value = -6350
$ws.Range("H43").Value = 1499.5
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("H113").Value = 7749.2856
$ws.Range("I113").Value = 6356.857
$ws.Range("J113").Value = 9141.714
$ws.Range("K113").Value = 6356.857
$ws.Range("L113").Value = 9141.714
$ws.Range("M113").Value = -3102.857
$ws.Range("N113").Value = -15649.714
$ws.Range("H132").Value = 999.5
$ws.Range("I132").Value = 999
$ws.Range("K132").Value = 2997
$ws.Range("M132").Value = -467
$ws.Range("H137").Value = 4693.385
$ws.Range("I137").Value = 4423.2
$ws.Range("J137").Value = 4862.25
$ws.Range("K137").Value = 13269.6
$ws.Range("L137").Value = 14586.75
$ws.Range("M137").Value = -10719.6
$ws.Range("N137").Value = -19686.75
$ws.Range("H141").Value = 2260.625
$ws.Range("I141").Value = 2298.1428
$ws.Range("K141").Value = 6894.428400000001
$ws.Range("M141").Value = -1714.428400000001

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 10808.454
$ws.Range("I2").Value = 923
$ws.Range("J2").Value = 37169.668
$ws.Range("K2").Value = 923
$ws.Range("L2").Value = 37169.668
$ws.Range("M2").Value = -810
$ws.Range("N2").Value = -37395.668
$ws.Range("H32").Value = 4577.9375
$ws.Range("I32").Value = 1603.3103
$ws.Range("J32").Value = 33332.668
$ws.Range("K32").Value = 1603.3103
$ws.Range("L32").Value = 33332.668
$ws.Range("M32").Value = -1316.3103
$ws.Range("N32").Value = -33906.668
$ws.Range("H45").Value = 874.8333
$ws.Range("I45").Value = 312.25
$ws.Range("J45").Value = 2000
$ws.Range("K45").Value = 312.25
$ws.Range("L45").Value = 2000
$ws.Range("M45").Value = 64.75
$ws.Range("N45").Value = -2754
$ws.Range("H61").Value = 3277.7144
$ws.Range("J61").Value = 3999.5
$ws.Range("L61").Value = 3999.5
$ws.Range("N61").Value = -4423.5
$ws.Range("H116").Value = 10808.454
$ws.Range("I116").Value = 923
$ws.Range("J116").Value = 37169.668
$ws.Range("K116").Value = 923
$ws.Range("L116").Value = 37169.668
$ws.Range("M116").Value = 1371
$ws.Range("N116").Value = -41757.668
$ws.Range("H132").Value = 2403.4285
$ws.Range("I132").Value = 2387.3333
$ws.Range("K132").Value = 7161.999899999999
$ws.Range("M132").Value = -4631.999899999999
$ws.Range("H136").Value = 3277.7144
$ws.Range("J136").Value = 3999.5
$ws.Range("L136").Value = 11998.5
$ws.Range("N136").Value = -17098.5

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 10808.454
$ws.Range("I3").Value = 923
$ws.Range("J3").Value = 37169.668
$ws.Range("K3").Value = 923
$ws.Range("L3").Value = 37169.668
$ws.Range("M3").Value = -809
$ws.Range("N3").Value = -37397.668
$ws.Range("H22").Value = 55
$ws.Range("I22").Value = 55
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 55
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 118
$ws.Range("N22").ClearContents()
$ws.Range("H99").Value = 1498.8572
$ws.Range("I99").Value = 1498.6666
$ws.Range("J99").Value = 1500
$ws.Range("K99").Value = 1498.6666
$ws.Range("L99").Value = 1500
$ws.Range("M99").Value = -0.6666000000000167
$ws.Range("N99").Value = -4496
$ws.Range("H139").Value = 45000
$ws.Range("I139").Value = 45000
$ws.Range("K139").Value = 45000
$ws.Range("M139").Value = -39860

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2254.88
$ws.Range("I31").Value = 1849.6923
$ws.Range("J31").Value = 2693.8333
$ws.Range("K31").Value = 1849.6923
$ws.Range("L31").Value = 2693.8333
$ws.Range("M31").Value = -1554.6923
$ws.Range("N31").Value = -3283.8333
$ws.Range("H34").Value = 2254.88
$ws.Range("I34").Value = 1849.6923
$ws.Range("J34").Value = 2693.8333
$ws.Range("K34").Value = 1849.6923
$ws.Range("L34").Value = 2693.8333
$ws.Range("M34").Value = -1647.6923
$ws.Range("N34").Value = -3097.8333
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H107").Value = 536.7273
$ws.Range("J107").Value = 682
$ws.Range("L107").Value = 682
$ws.Range("N107").Value = -4522
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()
$ws.Range("H134").Value = 7079.8
$ws.Range("I134").Value = 7079.8
$ws.Range("K134").Value = 21239.4
$ws.Range("M134").Value = -18704.4

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 10623.375
$ws.Range("I80").Value = 8999
$ws.Range("J80").Value = 10855.429
$ws.Range("K80").Value = 26997
$ws.Range("L80").Value = 32566.287
$ws.Range("M80").Value = -26061
$ws.Range("N80").Value = -34438.287
$ws.Range("H83").Value = 10623.375
$ws.Range("I83").Value = 8999
$ws.Range("J83").Value = 10855.429
$ws.Range("K83").Value = 80991
$ws.Range("L83").Value = 97698.861
$ws.Range("M83").Value = -76311
$ws.Range("N83").Value = -107058.861

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 369.8889
$ws.Range("J2").Value = 63.75
$ws.Range("L2").Value = 63.75
$ws.Range("N2").Value = -289.75
$ws.Range("H11").Value = 800012.6
$ws.Range("I11").Value = 400001.34
$ws.Range("J11").Value = 2600063.5
$ws.Range("K11").Value = 400001.34
$ws.Range("L11").Value = 2600063.5
$ws.Range("M11").Value = -399862.34
$ws.Range("N11").Value = -2600341.5
$ws.Range("H12").Value = 10685
$ws.Range("I12").Value = 10685
$ws.Range("K12").Value = 10685
$ws.Range("M12").Value = -10545
$ws.Range("H80").Value = 4000
$ws.Range("J80").Value = 4000
$ws.Range("L80").Value = 4000
$ws.Range("N80").Value = -5996
$ws.Range("H83").Value = 4000
$ws.Range("J83").Value = 4000
$ws.Range("L83").Value = 20000
$ws.Range("N83").Value = -29984
$ws.Range("H107").Value = 781.6667
$ws.Range("I107").Value = 781.6667
$ws.Range("K107").Value = 781.6667
$ws.Range("M107").Value = 1138.3333
$ws.Range("H122").Value = 35000
$ws.Range("J122").Value = 5000
$ws.Range("L122").Value = 15000
$ws.Range("N122").Value = -19900
$ws.Range("H132").Value = 2126.5
$ws.Range("I132").Value = 2126.5
$ws.Range("K132").Value = 6379.5
$ws.Range("M132").Value = -3849.5

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 23727.643
$ws.Range("I7").Value = 21476.334
$ws.Range("K7").Value = 21476.334
$ws.Range("M7").Value = -21364.334
$ws.Range("H100").Value = 7000
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 7000
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 7000
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -8082
$ws.Range("H103").Value = 27225
$ws.Range("J103").Value = 27225
$ws.Range("L103").Value = 27225
$ws.Range("N103").Value = -29569
$ws.Range("H126").Value = 23727.643
$ws.Range("I126").Value = 21476.334
$ws.Range("K126").Value = 64429.00199999999
$ws.Range("M126").Value = -61959.00199999999

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3200
$ws.Range("I62").Value = 3200
$ws.Range("K62").Value = 3200
$ws.Range("M62").Value = -2576
$ws.Range("H65").Value = 3200
$ws.Range("I65").Value = 3200
$ws.Range("K65").Value = 16000
$ws.Range("M65").Value = -12880
$ws.Range("H81").Value = 999
$ws.Range("I81").Value = 999
$ws.Range("K81").Value = 1998
$ws.Range("M81").Value = -937
$ws.Range("H84").Value = 999
$ws.Range("I84").Value = 999
$ws.Range("K84").Value = 9990
$ws.Range("M84").Value = -4686
$ws.Range("H113").Value = 203.83333
$ws.Range("I113").Value = 196.1
$ws.Range("J113").Value = 242.5
$ws.Range("K113").Value = 588.3
$ws.Range("L113").Value = 727.5
$ws.Range("M113").Value = 1581.7
$ws.Range("N113").Value = -5067.5
$ws.Range("H126").Value = 40166.555
$ws.Range("I126").Value = 37533.133
$ws.Range("K126").Value = 112599.399
$ws.Range("M126").Value = -110129.399
